# Update the stock-screener sheet: refresh ticker lists in columns B/C,
# clear stale entries in B/D/E/F, and extend the table down to row 24
# (column A keeps counting 0..22, using the same style as the existing
# rows for the new ones).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @{ A = index; B = ticker-or-$null; C = ticker-or-$null }
# D/E/F are cleared (set to empty) for every data row in the new range.
$rows = @(
    @{ R = 2;  A = 0;  B = "NSE:AKSHARCHEM"; C = "NSE:5PAISA" }
    @{ R = 3;  A = 1;  B = "NSE:BIGBLOC";    C = "NSE:AUTOBEES" }
    @{ R = 4;  A = 2;  B = "NSE:ETHOSLTD";   C = "NSE:BANCOINDIA" }
    @{ R = 5;  A = 3;  B = $null;            C = "NSE:BBL" }
    @{ R = 6;  A = 4;  B = $null;            C = "NSE:BERGEPAINT" }
    @{ R = 7;  A = 5;  B = $null;            C = "NSE:BIKAJI" }
    @{ R = 8;  A = 6;  B = $null;            C = "NSE:CHALET" }
    @{ R = 9;  A = 7;  B = $null;            C = "NSE:GLAND" }
    @{ R = 10; A = 8;  B = $null;            C = "NSE:GLOBUSSPR" }
    @{ R = 11; A = 9;  B = $null;            C = "NSE:GPIL" }
    @{ R = 12; A = 10; B = $null;            C = "NSE:HFCL" }
    @{ R = 13; A = 11; B = $null;            C = "NSE:IEX" }
    @{ R = 14; A = 12; B = $null;            C = "NSE:JAIPURKURT" }
    @{ R = 15; A = 13; B = $null;            C = "NSE:KICL" }
    @{ R = 16; A = 14; B = $null;            C = "NSE:KOTARISUG" }
    @{ R = 17; A = 15; B = $null;            C = "NSE:LINDEINDIA" }
    @{ R = 18; A = 16; B = $null;            C = "NSE:MID150BEES" }
    @{ R = 19; A = 17; B = $null;            C = "NSE:MSTCLTD" }
    @{ R = 20; A = 18; B = $null;            C = "NSE:MSUMI" }
    @{ R = 21; A = 19; B = $null;            C = "NSE:MURUDCERA" }
    @{ R = 22; A = 20; B = $null;            C = "NSE:PGHH" }
    @{ R = 23; A = 21; B = $null;            C = "NSE:RBLBANK" }
    @{ R = 24; A = 22; B = $null;            C = "NSE:RSWM" }
)

# Last row that already existed (and already carries the right
# formatting for column A) before this edit.
$lastExistingRow = 13

foreach ($row in $rows) {
    $r = $row.R

    $ws.Cells.Item($r, 1).Value = $row.A

    if ($r -gt $lastExistingRow) {
        # Brand-new row: copy column-A formatting (bold/border/centered
        # style) from the last pre-existing row so new rows match.
        $ws.Range("A$lastExistingRow").Copy()
        $ws.Range("A$r").PasteSpecial(-4122)
    }

    if ($row.B) {
        $ws.Cells.Item($r, 2).Value = $row.B
    } else {
        $ws.Cells.Item($r, 2).Value = ""
    }

    if ($row.C) {
        $ws.Cells.Item($r, 3).Value = $row.C
    } else {
        $ws.Cells.Item($r, 3).Value = ""
    }

    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = ""
}
